# Auto-generated script to update cryptos list (Price & Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "86.426.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +7.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.320.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "635.07"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.321"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +16.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.323.63"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.597"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000274"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.167"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.925.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.981.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.307.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.42"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +13.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.494.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.35"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000130"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +34.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "605.51"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.47"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.415"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.30"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +13.34%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.30%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "187.55"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.37"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.45%  "
